$d = $word.ActiveDocument

$replacements = @(
    @("2023-10-15 Sunday", "2023-10-16 Monday"),
    @("31×92=", "69×66="),
    @("50×89=", "86×44="),
    @("28×67=", "34×47="),
    @("59×63=", "29×56="),
    @("44×98=", "85×41="),
    @("49×52=", "79×99="),
    @("91×84=", "45×65="),
    @("79×61=", "86×28="),
    @("38×78=", "56×51="),
    @("91×87=", "59×83="),
    @("98×89=", "98×43="),
    @("66×59=", "65×28="),
    @("63×95=", "83×34="),
    @("47×72=", "34×23="),
    @("81×70=", "76×57="),
    @("16×83=", "37×69="),
    @("93×80=", "97×93="),
    @("46×99=", "33×33="),
    @("25×12=", "43×95="),
    @("21×72=", "44×33="),
    @("25×64=", "55×20="),
    @("93×19=", "61×94="),
    @("22×85=", "69×95="),
    @("45×44=", "11×23="),
    @("70×46=", "20×86=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
